$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after current row 8 (MOB) for "MOB CT" and "MOB PRE"
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# Copy formatting (border/font/alignment) from the MOB label cell (A8) onto the two new label cells
$ws.Cells.Item(8,1).Copy()
$ws.Cells.Item(9,1).PasteSpecial(-4122)
$ws.Cells.Item(10,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write out the refreshed activity data (rows 2-18)
$ws.Cells.Item(2,1).Value = "AMM"
$ws.Cells.Item(2,2).Value = 943
$ws.Cells.Item(2,3).Value = 25
$ws.Cells.Item(2,4).Value = 453.632
$ws.Cells.Item(2,5).Value = 1225.662
$ws.Cells.Item(2,6).Value = 481
$ws.Cells.Item(2,7).Value = 38
$ws.Cells.Item(2,8).Value = 14967
$ws.Cells.Item(2,9).Value = 9320
$ws.Cells.Item(2,10).Value = 38
$ws.Cells.Item(2,11).Value = 258
$ws.Cells.Item(2,12).Value = 2499

$ws.Cells.Item(3,1).Value = "AMM CT"
$ws.Cells.Item(3,2).Value = 271
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 0
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,11).Value = 0
$ws.Cells.Item(3,12).Value = 0

$ws.Cells.Item(4,1).Value = "COM"
$ws.Cells.Item(4,2).Value = 6
$ws.Cells.Item(4,3).Value = 11
$ws.Cells.Item(4,4).Value = 107
$ws.Cells.Item(4,5).Value = 211
$ws.Cells.Item(4,6).Value = 92
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 4381
$ws.Cells.Item(4,9).Value = 561
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(4,12).Value = 0

$ws.Cells.Item(5,1).Value = "IPR"
$ws.Cells.Item(5,2).Value = 1556
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = 131
$ws.Cells.Item(5,5).Value = 138
$ws.Cells.Item(5,6).Value = 7
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 2273
$ws.Cells.Item(5,9).Value = 26672
$ws.Cells.Item(5,10).Value = 2
$ws.Cells.Item(5,11).Value = 0
$ws.Cells.Item(5,12).Value = 0

$ws.Cells.Item(6,1).Value = "MIG"
$ws.Cells.Item(6,2).Value = 1895
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 464
$ws.Cells.Item(6,5).Value = 492
$ws.Cells.Item(6,6).Value = 28
$ws.Cells.Item(6,7).Value = 4
$ws.Cells.Item(6,8).Value = 2867
$ws.Cells.Item(6,9).Value = 36871
$ws.Cells.Item(6,10).Value = 17
$ws.Cells.Item(6,11).Value = 0
$ws.Cells.Item(6,12).Value = 0

$ws.Cells.Item(7,1).Value = "MKT"
$ws.Cells.Item(7,2).Value = 140
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = 78
$ws.Cells.Item(7,5).Value = 84
$ws.Cells.Item(7,6).Value = 6
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = 714
$ws.Cells.Item(7,9).Value = 6184
$ws.Cells.Item(7,10).Value = 2
$ws.Cells.Item(7,11).Value = 0
$ws.Cells.Item(7,12).Value = 0

$ws.Cells.Item(8,1).Value = "MOB"
$ws.Cells.Item(8,2).Value = 809
$ws.Cells.Item(8,3).Value = 11
$ws.Cells.Item(8,4).Value = 767
$ws.Cells.Item(8,5).Value = 1189
$ws.Cells.Item(8,6).Value = 341
$ws.Cells.Item(8,7).Value = 12
$ws.Cells.Item(8,8).Value = 4381
$ws.Cells.Item(8,9).Value = 7052
$ws.Cells.Item(8,10).Value = 13
$ws.Cells.Item(8,11).Value = 58
$ws.Cells.Item(8,12).Value = 731

$ws.Cells.Item(9,1).Value = "MOB CT"
$ws.Cells.Item(9,2).Value = 270
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 4
$ws.Cells.Item(9,5).Value = 4
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(9,8).Value = 0
$ws.Cells.Item(9,9).Value = 2500
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = 0
$ws.Cells.Item(9,12).Value = 0

$ws.Cells.Item(10,1).Value = "MOB PRE"
$ws.Cells.Item(10,2).Value = 1048
$ws.Cells.Item(10,3).Value = 8
$ws.Cells.Item(10,4).Value = 744
$ws.Cells.Item(10,5).Value = 806
$ws.Cells.Item(10,6).Value = 49
$ws.Cells.Item(10,7).Value = 6
$ws.Cells.Item(10,8).Value = 6462
$ws.Cells.Item(10,9).Value = 13839
$ws.Cells.Item(10,10).Value = 7
$ws.Cells.Item(10,11).Value = 0
$ws.Cells.Item(10,12).Value = 0

$ws.Cells.Item(11,1).Value = "MSK"
$ws.Cells.Item(11,2).Value = 1056
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 275
$ws.Cells.Item(11,5).Value = 299
$ws.Cells.Item(11,6).Value = 24
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 968
$ws.Cells.Item(11,9).Value = 10200
$ws.Cells.Item(11,10).Value = 1
$ws.Cells.Item(11,11).Value = 0
$ws.Cells.Item(11,12).Value = 0

$ws.Cells.Item(12,1).Value = "NOT"
$ws.Cells.Item(12,2).Value = 1405
$ws.Cells.Item(12,3).Value = 0
$ws.Cells.Item(12,4).Value = 82
$ws.Cells.Item(12,5).Value = 93
$ws.Cells.Item(12,6).Value = 11
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = 4367
$ws.Cells.Item(12,9).Value = 6603
$ws.Cells.Item(12,10).Value = 2
$ws.Cells.Item(12,11).Value = 0
$ws.Cells.Item(12,12).Value = 0

$ws.Cells.Item(13,1).Value = "TEC"
$ws.Cells.Item(13,2).Value = 826
$ws.Cells.Item(13,3).Value = 1
$ws.Cells.Item(13,4).Value = 502
$ws.Cells.Item(13,5).Value = 529
$ws.Cells.Item(13,6).Value = 22
$ws.Cells.Item(13,7).Value = 4
$ws.Cells.Item(13,8).Value = 3968
$ws.Cells.Item(13,9).Value = 6224
$ws.Cells.Item(13,10).Value = 7
$ws.Cells.Item(13,11).Value = 0
$ws.Cells.Item(13,12).Value = 0

$ws.Cells.Item(14,1).Value = "TEC CT"
$ws.Cells.Item(14,2).Value = 272
$ws.Cells.Item(14,3).Value = 0
$ws.Cells.Item(14,4).Value = 20
$ws.Cells.Item(14,5).Value = 24
$ws.Cells.Item(14,6).Value = 4
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = 1667
$ws.Cells.Item(14,9).Value = 3000
$ws.Cells.Item(14,10).Value = 0
$ws.Cells.Item(14,11).Value = 0
$ws.Cells.Item(14,12).Value = 0

$ws.Cells.Item(15,1).Value = "TST"
$ws.Cells.Item(15,2).Value = 502
$ws.Cells.Item(15,3).Value = 0
$ws.Cells.Item(15,4).Value = 70
$ws.Cells.Item(15,5).Value = 95
$ws.Cells.Item(15,6).Value = 25
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = 8404
$ws.Cells.Item(15,9).Value = 25094
$ws.Cells.Item(15,10).Value = 1
$ws.Cells.Item(15,11).Value = 0
$ws.Cells.Item(15,12).Value = 0

$ws.Cells.Item(16,1).Value = "VIP"
$ws.Cells.Item(16,2).Value = 194
$ws.Cells.Item(16,3).Value = 0
$ws.Cells.Item(16,4).Value = 2
$ws.Cells.Item(16,5).Value = 4
$ws.Cells.Item(16,6).Value = 2
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = 50
$ws.Cells.Item(16,9).Value = 0
$ws.Cells.Item(16,10).Value = 0
$ws.Cells.Item(16,11).Value = 0
$ws.Cells.Item(16,12).Value = 0

$ws.Cells.Item(17,1).Value = "WLC"
$ws.Cells.Item(17,2).Value = 1823
$ws.Cells.Item(17,3).Value = 0
$ws.Cells.Item(17,4).Value = 81
$ws.Cells.Item(17,5).Value = 87
$ws.Cells.Item(17,6).Value = 6
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = 2857
$ws.Cells.Item(17,9).Value = 17004
$ws.Cells.Item(17,10).Value = 0
$ws.Cells.Item(17,11).Value = 0
$ws.Cells.Item(17,12).Value = 0

$ws.Cells.Item(18,1).Value = "ZERO"
$ws.Cells.Item(18,2).Value = 3673
$ws.Cells.Item(18,3).Value = 51
$ws.Cells.Item(18,4).Value = 561
$ws.Cells.Item(18,5).Value = 2482
$ws.Cells.Item(18,6).Value = 943
$ws.Cells.Item(18,7).Value = 27
$ws.Cells.Item(18,8).Value = 14282
$ws.Cells.Item(18,9).Value = 15925
$ws.Cells.Item(18,10).Value = 3
$ws.Cells.Item(18,11).Value = 900
$ws.Cells.Item(18,12).Value = 16154

